# Build new site after additions
# Applies the cell-value additions recorded in the commit diff: a handful of
# 0 -> 1 flags (and one 1 -> 0) across the "Tabelle1" matrix, plus moving the
# active selection to L39 to match the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 2
$ws.Range("L2").Value = 1

# Row 4
$ws.Range("L4").Value = 1
$ws.Range("N4").Value = 1

# Row 10
$ws.Range("N10").Value = 1

# Row 12
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1

# Row 15
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 1

# Row 16
$ws.Range("K16").Value = 1

# Row 17
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 1

# Row 18
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 1

# Row 19
$ws.Range("L19").Value = 1
$ws.Range("N19").Value = 1

# Row 20
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 1

# Row 21
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 1

# Row 22
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("N22").Value = 1

# Row 23
$ws.Range("L23").Value = 1
$ws.Range("N23").Value = 1

# Row 24
$ws.Range("N24").Value = 0

# Row 25
$ws.Range("L25").Value = 1

# Row 26
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 1

# Row 27
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 1

# Row 32
$ws.Range("L32").Value = 1

# Row 33
$ws.Range("N33").Value = 1

# Row 36
$ws.Range("N36").Value = 1

# Row 37
$ws.Range("L37").Value = 1
$ws.Range("N37").Value = 1

# Row 38
$ws.Range("K38").Value = 1
$ws.Range("L38").Value = 1
$ws.Range("N38").Value = 1

# Row 40
$ws.Range("K40").Value = 1

# Row 41
$ws.Range("K41").Value = 1
$ws.Range("N41").Value = 1

# Move the active selection/cursor to match the saved workbook state.
$ws.Range("L39").Select()
